$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new
# header cells so they pick up the same style (bold, border, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the data values for columns I (I0) and J (IF), rows 2-10
$dataI = @(8, 8, 9, 10, 6, 8, 1, 3, 3)
$dataJ = @(9, 9, 9, 10, 7, 8, 3, 3, 4)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
